# Insert a new data row at row 436 (pushes existing rows 436-468 down to 437-469)
# and populate it with a new record, matching the rest of the table's columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 436.
$ws.Rows.Item(436).Insert()

# Populate the new row 436 with the same field layout as the surrounding rows,
# using the new date (2023-08-28 -> serial 45166) and volume (1400).
$ws.Range("A436").Value = 8
$ws.Range("B436").Value = "Terminal La Palmera de La Serena"
$ws.Range("C436").Value = "Coquimbo"
$ws.Range("D436").Value = 45166
$ws.Range("E436").Value = 4
$ws.Range("F436").Value = 100112012
$ws.Range("G436").Value = "Espinaca"
$ws.Range("H436").Value = "Sin especificar"
$ws.Range("I436").Value = "Primera"
$ws.Range("J436").Value = 1400
$ws.Range("K436").Value = 450
$ws.Range("L436").Value = 500
$ws.Range("M436").Value = 475
$ws.Range("N436").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O436").Value = "Provincia del Elquí"
$ws.Range("P436").Value = 950
$ws.Range("Q436").Value = 0.5
$ws.Range("R436").Value = "Hortaliza"
